# Auto-generated edit script: refresh market-board derived values
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) per the scheduled-runner sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 11000
$ws.Range("J113").Value = 11000
$ws.Range("L113").Value = 11000
$ws.Range("N113").Value = -17508

$ws.Range("H132").Value = 55561004
$ws.Range("I132").Value = 55561004
$ws.Range("K132").Value = 166683012
$ws.Range("M132").Value = -166680482

$ws.Range("H137").Value = 17545572
$ws.Range("I137").Value = 27778864
$ws.Range("J137").Value = 2785.5715
$ws.Range("K137").Value = 83336592
$ws.Range("L137").Value = 8356.7145
$ws.Range("M137").Value = -83334042
$ws.Range("N137").Value = -13456.7145

$ws.Range("H138").Value = 3764.5757
$ws.Range("I138").Value = 4369.077
$ws.Range("J138").Value = 3371.65
$ws.Range("K138").Value = 13107.231
$ws.Range("L138").Value = 10114.95
$ws.Range("M138").Value = -7967.231
$ws.Range("N138").Value = -20394.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2255.3076
$ws.Range("I2").Value = 1646.6666
$ws.Range("K2").Value = 1646.6666
$ws.Range("M2").Value = -1533.6666

$ws.Range("H32").Value = 17133.428
$ws.Range("I32").Value = 20032.928
$ws.Range("K32").Value = 20032.928
$ws.Range("M32").Value = -19745.928

$ws.Range("H61").Value = 6164.4
$ws.Range("I61").Value = 8374.333000000001
$ws.Range("K61").Value = 8374.333000000001
$ws.Range("M61").Value = -8162.333000000001

$ws.Range("H110").Value = 1242.5834
$ws.Range("I110").Value = 1242.5834
$ws.Range("K110").Value = 1242.5834
$ws.Range("M110").Value = 802.4166

$ws.Range("H116").Value = 2255.3076
$ws.Range("I116").Value = 1646.6666
$ws.Range("K116").Value = 1646.6666
$ws.Range("M116").Value = 647.3334

$ws.Range("H131").Value = 79996.8
$ws.Range("J131").Value = 79996.8
$ws.Range("L131").Value = 79996.8
$ws.Range("N131").Value = -90076.8

$ws.Range("H132").Value = 4818.615
$ws.Range("J132").Value = 2193.5
$ws.Range("L132").Value = 6580.5
$ws.Range("N132").Value = -11640.5

$ws.Range("H136").Value = 6164.4
$ws.Range("I136").Value = 8374.333000000001
$ws.Range("K136").Value = 25122.999
$ws.Range("M136").Value = -22572.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2255.3076
$ws.Range("I3").Value = 1646.6666
$ws.Range("K3").Value = 1646.6666
$ws.Range("M3").Value = -1532.6666

$ws.Range("H86").Value = 18185.334
$ws.Range("I86").Value = 16620.916
$ws.Range("J86").Value = 20688.4
$ws.Range("K86").Value = 16620.916
$ws.Range("L86").Value = 20688.4
$ws.Range("M86").Value = -15497.916
$ws.Range("N86").Value = -22934.4

$ws.Range("H89").Value = 18185.334
$ws.Range("I89").Value = 16620.916
$ws.Range("J89").Value = 20688.4
$ws.Range("K89").Value = 83104.58
$ws.Range("L89").Value = 103442
$ws.Range("M89").Value = -77488.58
$ws.Range("N89").Value = -114674

$ws.Range("H99").Value = 1942.5
$ws.Range("I99").Value = 1942.5
$ws.Range("K99").Value = 1942.5
$ws.Range("M99").Value = -444.5

$ws.Range("H105").Value = 2727.4285
$ws.Range("I105").Value = 2727.4285
$ws.Range("K105").Value = 2727.4285
$ws.Range("M105").Value = -980.4285

$ws.Range("H134").Value = 996.875
$ws.Range("I134").Value = 996.875
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2990.625
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -455.625
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("K16").Value = 2000
$ws.Range("M16").Value = -1713

$ws.Range("H99").Value = 3647.9092
$ws.Range("J99").Value = 3900.7144
$ws.Range("L99").Value = 3900.7144
$ws.Range("N99").Value = -6896.7144

$ws.Range("H107").Value = 2068.2856
$ws.Range("I107").Value = 1996.3334
$ws.Range("K107").Value = 1996.3334
$ws.Range("M107").Value = -76.33339999999998

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

$ws.Range("H126").Value = 3647.9092
$ws.Range("J126").Value = 3900.7144
$ws.Range("L126").Value = 11702.1432
$ws.Range("N126").Value = -16642.1432

$ws.Range("H132").Value = 40023016
$ws.Range("I132").Value = 50028296
$ws.Range("J132").Value = 1897
$ws.Range("K132").Value = 150084888
$ws.Range("L132").Value = 5691
$ws.Range("M132").Value = -150082358
$ws.Range("N132").Value = -10751

$ws.Range("H134").Value = 1839.0344
$ws.Range("I134").Value = 1939.45
$ws.Range("K134").Value = 5818.35
$ws.Range("M134").Value = -3283.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 492.8889
$ws.Range("I6").Value = 54.5
$ws.Range("K6").Value = 163.5
$ws.Range("M6").Value = -50.5

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H99").Value = 2159.1785
$ws.Range("I99").Value = 1328.5385
$ws.Range("J99").Value = 2879.0667
$ws.Range("K99").Value = 3985.6155
$ws.Range("L99").Value = 8637.2001
$ws.Range("M99").Value = -1739.6155
$ws.Range("N99").Value = -13129.2001

$ws.Range("H113").Value = 2506.9285
$ws.Range("J113").Value = 2627.4285
$ws.Range("L113").Value = 7882.2855
$ws.Range("N113").Value = -12222.2855

$ws.Range("H131").Value = 1975.4445
$ws.Range("J131").Value = 4748.5
$ws.Range("L131").Value = 14245.5
$ws.Range("N131").Value = -24325.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17548332
$ws.Range("I132").Value = 4775.6875
$ws.Range("K132").Value = 14327.0625
$ws.Range("M132").Value = -11797.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41670620
$ws.Range("I40").Value = 71431150
$ws.Range("J40").Value = 5872.6
$ws.Range("K40").Value = 71431150
$ws.Range("L40").Value = 5872.6
$ws.Range("M40").Value = -71431014
$ws.Range("N40").Value = -6144.6

$ws.Range("H43").Value = 44999
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H61").Value = 3169.75
$ws.Range("I61").Value = 2810.2632
$ws.Range("K61").Value = 2810.2632
$ws.Range("M61").Value = -2608.2632

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H113").Value = 3169.75
$ws.Range("I113").Value = 2810.2632
$ws.Range("K113").Value = 2810.2632
$ws.Range("M113").Value = -640.2631999999999

$ws.Range("H122").Value = 3362.4285
$ws.Range("I122").Value = 3089.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9268.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6818.5
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 5356.7144
$ws.Range("J132").Value = 5666.3335
$ws.Range("L132").Value = 16999.0005
$ws.Range("N132").Value = -22059.0005

$ws.Range("H136").Value = 5000
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 14000
$ws.Range("I38").Value = 8000
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 8000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = -7527
$ws.Range("N38").Value = -20946

$ws.Range("H95").Value = 19584.5
$ws.Range("J95").Value = 19584.5
$ws.Range("L95").Value = 19584.5
$ws.Range("N95").Value = -25076.5

$ws.Range("H96").Value = 1963.125
$ws.Range("I96").Value = 2027.4
$ws.Range("J96").Value = 999
$ws.Range("K96").Value = 2027.4
$ws.Range("L96").Value = 999
$ws.Range("M96").Value = -654.4000000000001
$ws.Range("N96").Value = -3745

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 41676760
$ws.Range("I132").Value = 12239.263
$ws.Range("K132").Value = 36717.789
$ws.Range("M132").Value = -34187.789
